# Contest schedule update: populate match rows 22 to 35 (sheet rows 34-47)
# on Sheet1 with the "Format" team-vs-team labels, a Points flag of 1, and
# the same VLOOKUP/RANK "Rank" formulas used by the rows above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row (sheet row number) -> match label, in contest order 22..35
$matches = @{
    34 = "CSK vs RCB"
    35 = "MI vs PBKS"
    36 = "RR vs GT"
    37 = "SRH vs KKR"
    38 = "MI vs LSG"
    39 = "DC vs RCB"
    40 = "PBKS vs SRH"
    41 = "GT vs CSK"
    42 = "RR vs KKR"
    43 = "LSG vs RCB"
    44 = "DC vs PBKS"
    45 = "MI vs CSK"
    46 = "DC vs RR"
    47 = "KKR vs GT"
}

# Columns whose formula cell (left) ranks the score in the column right after it.
$scoreCols = @("E", "H", "K", "N", "Q", "T", "W")
$rankCols  = @("D", "G", "J", "M", "P", "S", "V")

foreach ($r in 34..47) {
    $ws.Range("B$r").Value = 1
    $ws.Range("C$r").Value = $matches[$r]

    for ($i = 0; $i -lt $scoreCols.Length; $i++) {
        $scoreCol = $scoreCols[$i]
        $rankCol = $rankCols[$i]
        $formula = '=IF(ISERROR(VLOOKUP(RANK(' + $scoreCol + $r + ', ($W' + $r + ',$T' + $r + ',$Q' + $r + ',$N' + $r + ',$K' + $r + ',$H' + $r + ',$E' + $r + '), 0),  $A$2:$B$10, 2, FALSE)),"",VLOOKUP(RANK(' + $scoreCol + $r + ', ($W' + $r + ',$T' + $r + ',$Q' + $r + ',$N' + $r + ',$K' + $r + ',$H' + $r + ',$E' + $r + '), 0),  $A$2:$B$10, 2, FALSE))'
        $ws.Range($rankCol + $r).Formula = $formula
    }
}
